$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "My Series" to "Data"
$ws.Name = "Data"

# Update the refreshed-series values in column A (rows 13-22)
$ws.Range("A13").Value = 15733.82555555556
$ws.Range("A14").Value = 139798646.1781028
$ws.Range("A15").Value = 11823.64775262282
$ws.Range("A16").Value = 0.951248095277063
$ws.Range("A17").Value = -0.1855642624124751
$ws.Range("A18").Value = 0.7514795248538862
$ws.Range("A19").Value = 3271.34
$ws.Range("A21").Value = 10540.67
$ws.Range("A22").Value = 9

# Row 23 keeps the "###0.000" number format but its value changes too
$ws.Range("A23").Value = 3271.34
$ws.Range("A23").NumberFormat = "###0.000"

# New data points appended below the previous last row
$ws.Range("A24").Value = 6384.48
$ws.Range("A24").NumberFormat = "###0.000"
$ws.Range("A25").Value = 7111.74
$ws.Range("A25").NumberFormat = "###0.000"
$ws.Range("A26").Value = 8354.02
$ws.Range("A26").NumberFormat = "###0.000"
$ws.Range("A27").Value = 10540.67
$ws.Range("A27").NumberFormat = "###0.000"
$ws.Range("A28").Value = 16421.36
$ws.Range("A28").NumberFormat = "###0.000"
$ws.Range("A29").Value = 22352.63
$ws.Range("A29").NumberFormat = "###0.000"
$ws.Range("A30").Value = 29205.22
$ws.Range("A30").NumberFormat = "###0.000"
$ws.Range("A31").Value = 37962.97
$ws.Range("A31").NumberFormat = "###0.000"

# Refresh the embedded data-provider comment payload on A1 (opaque vendor blob)
$c = $ws.Comments.Item(1)
$c.Text("QxwAAB+LCAAAAAAAAAOlGdtuI0n1V1p+Agm72/ZkclGlV75l1sKOI9shm31B5e5K3KS72nRVJ/HbIoEWLYsQQrNouYqnRUgMI9iVlhku/7KaZIYnfoFTl77Zzk46jKJJ17nVqVPnWkHvXAe+cUki5oV0v1KvWRWDUCd0PXq+X4n5WbX+uPKOjXrXDvGPcIQDwoHYAC7K9q6Zt1+Zc77YM82rq6vaVbMWRudmw7Lq5nvDwcSZkwBXPco4pg6ppFzu27kqNuq4wZBw7GKOFed+pT/p1zrEc7oAG2KKz0lUa8fMo4SxHuUe9wgTnBHBnHS6w++og9mN2uNaHZlr8IyyHXu+q+gKlAqu6WBbMvUCYjes+k7V2q1azam1tdds7DW3a7u7jfcTxpQQDTDjExJdeo4ETDgOFpLd2rWa1laz0dxG5kYikJUZwEYj3x2TS48Rt0N8n5WyiKkvsOVwOHU5Y1rIzPFqQQ9X4UmEF/Opx31STo3xsG0EVOuSCbHRQRgRB+z3IJUOydUo0madLgaAnc69iC+7eFla1jEj0WghjFSO1UbdkPKWTyJ+vIC7Ji64AiBsHsUEmXcgM6auxxz49mhMXPsM+yzPVECikzC6YAvskEOIY1PIuKJ+iF1wOO4x7jnZpmsIdBSFC5AIm7dD3z0AqZp4AyKV3KdgYrFtOwwvMu02IZG8VXm/cKcB5gn5GhxN5uHViPrLSTxjTuTNiNttJ9QbcUgEpObuxIyHAWiRgZCC5SBL+AcBuApGXeJ4AfaPfDAis5sgpQBArZiHZx7vhH4cUJbotAJFJ3CiKblOT5iu0Qgulwqjh7RPE3pl5o2oIsM4vEr3XEdII+TALeYk172OWCXuAiy5vnWMvBFxygPPh/qQv4sctOgVkzkhfKNLKAwSqfBAVBy7vRR7IjODIPBM8G6A2vXd3e2qVYefqWXtyR/YOUWjHnXlB9SXx9V6o9rM0yVIdBgHoxlE8KU8k10H3AoIwSn8to/pBUBPPD4/bCXab8AgdeY76ddxCGJ14eOlBKd2ycNQnzp+7BKVAvr0TDql0E1d451otAYaQFTbCNPldLmATMy8PQ4f+xWozXuMR1D9K7YTxpRHS5ErkKlJ38bD4hmVG2D/3jxnEfl+DE3H8iCmTid077+bq6xzTD1+fw3DOFIJ8P4s0noiF8asS0RWkWn+3vxOmTOxqBR5QEkQUs+5v7XByEJ79wEHYUlU5TjMO3yrS85w7ENbw6H6nKdpaRWMWuxilSYPQseRnyQHWzSNDLpGxw1qDtRV0RnVnDAQABOatZMJMvP0ojlwSI+eDzA9j6H8pgG4Ck9Tkygd0whTJo6TVtuVLLWZCCUBrboAW0X5KJYWU1EeAhaZK3RoSoJFGGF/CIbxDvT96FYCqvQQ87leQdr3iZMY2cxYU66iZonibyOT+VsdQ0SGzicrQEkkzqL601yOz4BIHHMIDux3sO/NIpV/kjK3CQc3ljVOSaYSpyvZRCWXAEMKVKZvk6XoXLOFhkufrScI5cEi5diT8aOdxpbVbEDVF2skjzwm2Dd64PacGH16SRgPgG3PGBPmufDlYX/PeJfMiAflQtpIJ+zS3Hk+dJBkRKlKC/QtQooEUIvPPUi464QpJmOwTwmO/GWOUB11EDpAd/vjf9/85vmrF5/dfvT0zRc//O8/fvXqnz+7efYj+Lj9699uPv6lOqYiRlM884lUaNre2bGaj8DRUhASxjVl2+jGDpew01PZLaZrpIceuej0+p0ng7ZMKCkwYVfJ1xTz1DKMs+VEHUJuJK/UTDxBkdjTJEHpdQGbS+a2mG8uSZE6j7+LUdni9cvPXr/8853c2mCFpmVLNCNvb1rq9TW6tGkRc6MaC5Kp9FHV2qo2GjniFRo0huEapoDUTn3XbtZhGG00rXqazN3UkTcRraK0pCk+N1f4FKijGonUBfLrBCkdfwohkqJVKOQW2kW/+MmbvzwtUGnrakhRCignC77YzEwWUvTheGpMRsfjTs+Y9ibCTzJcjk4J/xpivXsaTwWnojTG/rcMQsWDilGBaaFihGcGwc7cWEIk5uKw4GyboGqjB4pc1fJJFMYLdSM5hgy6gTLNJhs5NuQaiZP2XEs6GWoDudL15u+fb2LQB+lmrV/6NpCHoQJGgXJ4HbWf/uvVlx++evHi9vnPb778QUGC3icdmcHPIZryy9TtIeXperMCQScTacwL67u5+qKBYuQ4Cj3KmV1/LKcNvULAWhfS5G/UD6DkScHSXgBfgaB3Metdcx3Y9iEyiwDQc4Gh2obZXJYCVA7P7Pqf3/7u9tef337y/M2Hf7r56I83H3/y+uXv3zz7g4q626fPb3/6TGf51UIgdRHTnuoCDfl24BgiGg1Ru42vPviFQUNuQM9hxDIjffXBpzlhQlHZnWSSoadLFSmqsEaaZxZ8Rk6VVIcCX8qiGoCOKGHNlEIXsXDhOdkm71eFKBF3EvGN/rQaM2KE0E59E05SJM6Y78unWVRJPdq2GvWGxiptxBFmmOVM/8QPZ9BkJAg5nK+QFLi+niGjlfs9GYzarUFGopQYRS7M8pZ4YRAfKOkpRUnps2SVuFoOAlho/JzYF+8pa2TrqFRyLo2Z+mHirOWK9Gdblpz8i/WxQIE6cRSphojqd+5JvIBuOHm+uhsvn/RyDfCh6lXzLXG27neLeFjnsFAIi2gBkHiZmjRKpak+E28gqp09FKbJloArPAOCOfRTtuq0LqGvjEyRd3pRFEYbk0+GSciG0ElDRjEzi6c08k5V1+1md5UAkoSXfqjRT58w7BKf8HLvvGbGPQwvH8wLd1+Wtc9GvquNWW70SM2SCcg/dgtH+X/fupWztaIIGivxOFb6cTqZXMcw8JbURh1FMooREHbXr8cHXsT4eyIT6C8FOU0hp6pDlSTqQ65P7aYGAIGZl24W1ExCl6s/OYT+wAu8kmOhlcR3UQjYcrFQLVy/nKeI0nJIrqHBzEmApDj7HpQN9XRSRppyWMilKb945WPe+ZyXVWx7holLZlbVmZFG9ZFr7VR3CWlW63X4HzuNhmVtiTdCLRwyh0euSm5iJheW/SnQ/h+Xl9fgQxwAAA==")
